$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column A with sequential numbers 1-7 for rows 2-8
for ($i = 0; $i -le 6; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $i + 1
}

# Update the selection to A8 (matches the new <selection activeCell="A8" sqref="A8"/>)
$ws.Range("A8").Select()

# Best-effort: try to nudge the recorded window size to match the saved
# workbook view geometry (harmless no-op if unsupported by the host).
try { $excel.ActiveWindow.Height = 9287 } catch {}
